# Update the cryptos list sheet with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number need to be forced to stay
# text (matching the source workbook, where every Price/Volume cell is
# stored as a string) - otherwise Excel's COM layer auto-converts them.
$textCells = @(
  "D5","D6","D7","D10","D11","D12","D15","D17","D21","D22","D23",
  "D26","D27","D28","D29","D30","D32","D36","D37","D38","D39","D41",
  "D42","D43","D44","D45","D46","D47","D48","D50","D51"
)
foreach ($addr in $textCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.919.24"
$ws.Range("E2").Value = "  +0.52%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.262.43"
$ws.Range("E3").Value = "  -0.62%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - now BNB (was XRP)
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "233.15"
$ws.Range("E5").Value = "  +0.82%  "

# Row 6 - now XRP (was BNB)
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "0.652"
$ws.Range("E6").Value = "  +3.85%  "

# Row 7 - Solana
$ws.Range("D7").Value = "63.82"
$ws.Range("E7").Value = "  +0.11%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.08%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +4.91%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0974"
$ws.Range("E10").Value = "  -6.98%  "

# Row 11 - OKB
$ws.Range("D11").Value = "58.28"
$ws.Range("E11").Value = "  +1.70%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "26.47"
$ws.Range("E12").Value = "  +1.83%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.91%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.598.35"
$ws.Range("E14").Value = "  -0.60%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "15.64"

# Row 16 - Polkadot
$ws.Range("E16").Value = "  +4.01%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.844"
$ws.Range("E17").Value = "  +2.43%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.264.46"
$ws.Range("E18").Value = "  +0.04%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "43.820.63"
$ws.Range("E19").Value = "  +0.33%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0980"
$ws.Range("E20").Value = "  -2.76%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "73.82"
$ws.Range("E21").Value = "  +0.46%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.18"
$ws.Range("E22").Value = "  +1.41%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "249.93"
$ws.Range("E23").Value = "  -0.27%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.00%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -1.52%  "

# Row 26 - now WEMIXToken (was Toncoin)
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "3.56"
$ws.Range("E26").Value = "  +27.33%  "

# Row 27 - now Toncoin (was WEMIXToken)
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "2.24"
$ws.Range("E27").Value = "  -3.44%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "9.92"
$ws.Range("E28").Value = "  +0.56%  "

# Row 29 - Monero
$ws.Range("D29").Value = "173.97"

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "21.93"
$ws.Range("E30").Value = "  +4.68%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  +0.19%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "1.44"
$ws.Range("E32").Value = "  -0.32%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  +3.04%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  +4.62%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  +0.05%  "

# Row 36 - InternetComputer(DFINITY)
$ws.Range("D36").Value = "4.96"
$ws.Range("E36").Value = "  -1.85%  "

# Row 37 - RenderToken
$ws.Range("D37").Value = "3.71"
$ws.Range("E37").Value = "  -2.79%  "

# Row 38 - THORChain
$ws.Range("D38").Value = "6.43"
$ws.Range("E38").Value = "  -5.19%  "

# Row 39 - LidoDAOToken
$ws.Range("D39").Value = "2.30"
$ws.Range("E39").Value = "  -1.82%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +3.11%  "

# Row 41 - BinanceUSD
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.05%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "8.77"
$ws.Range("E42").Value = "  +4.55%  "

# Row 43 - now FTXToken (was InjectiveProtocol)
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").Value = "4.53"
$ws.Range("E43").Value = "  +2.47%  "

# Row 44 - now InjectiveProtocol (was Aave)
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "17.38"
$ws.Range("E44").Value = "  +0.95%  "

# Row 45 - now Aave (was Cronos)
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "98.83"
$ws.Range("E45").Value = "  +1.23%  "

# Row 46 - TrustWalletToken
$ws.Range("D46").Value = "1.20"
$ws.Range("E46").Value = "  -0.76%  "

# Row 47 - now Cronos (was FTXToken)
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.0953"
$ws.Range("E47").Value = "  -0.90%  "

# Row 48 - TerraClassic
$ws.Range("D48").Value = "0.000211"
$ws.Range("E48").Value = "  +2.52%  "

# Row 49 - Maker
$ws.Range("D49").Value = "1.458.91"
$ws.Range("E49").Value = "  -1.36%  "

# Row 50 - NEARProtocol
$ws.Range("D50").Value = "2.32"
$ws.Range("E50").Value = "  -0.80%  "

# Row 51 - Celestia
$ws.Range("D51").Value = "9.99"
$ws.Range("E51").Value = "  -4.68%  "
